$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "insolvencni rizeni: soudni rizeni" -> add "(insolvence)"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "insolvenční řízení: soudní řízení",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "insolvenční řízení (insolvence): soudní řízení", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Aby mohlo zapocit insolvencni rizeni rozhoduje soud..." -> reworded
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Aby mohlo započít insolvenční řízení rozhoduje soud, zda je dlužník v",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Aby mohl soud rozhodnout o způsobu úpadku v insolvenční řízení rozhoduje soud, zda je dlužník v", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) ", byt v upadku je tedy zakladni prerekvizitou ... (insolvence)." -> reworded ending
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    ", být v úpadku je tedy základní prerekvizitou pro to, aby mohlo proběhnout insolvenční řízení (insolvence).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", být v úpadku je tedy základní prerekvizitou pro to, aby mohl soud rozhodnout o způsobu úpadku.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "- vyhlasuje se rozhodnutim:" -> extended explanation
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "- vyhlašuje se rozhodnutím:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- vyhlašuje se rozhodnutím ve formě usnesení, přičemž opravné prostředky nejsou standardně přípustné, pokud zákon nestanoví jinak:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) "zakaz jednani vedouciho ... zakon." -> ends with comma, + a new bullet after it
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "zákaz jednání vedoucího k uspokojení pohledávek mimo insolvenční zákon, ledaže to povoluje zákon.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "zákaz jednání vedoucího k uspokojení pohledávek mimo insolvenční zákon, ledaže to povoluje zákon,", 2) | Out-Null

$rngZakaz = $d.Content
$rngZakaz.Find.Execute(
    "zákaz jednání vedoucího k uspokojení pohledávek mimo insolvenční zákon, ledaže to povoluje zákon,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$rngZakaz.Collapse(0)
$rngZakaz.InsertAfter("`r" + "jednotnost - další insolvenční návrh - přistoupení k řízení.")

# ---------------------------------------------------------------------------
# 6) append the new "Prubeh insolvence" section (with subsections) at the
#    very end of the document body, including headings, bookmarks and the
#    new numbered lists (numId 1014 / 1015 / 1016).
# ---------------------------------------------------------------------------
$bodyXml = @'
<w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:bookmarkStart w:id="43" w:name="průběh-insolvence"/><w:bookmarkEnd w:id="43"/><w:r><w:t xml:space="preserve">Průběh insolvence</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1014"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Lze zahájit pouze na návrh, zahajuje se dnem kdy dojde návrh na zahájení insolvenčního řízení pokud ho podává dlužník a chce dosáhnout odlužení, musí rovněž připojit návrh na oddlužení.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1014"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Návrh může podat dlužník, nebo věřitel, jde-li o hrozící úpadek, může jej podat pouze dlužník.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1014"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Dlužník (popřípadě zákonní zástupci dlužníka, jeho statutární orgán a likvidátor), který je právnickou osobou nebo fyzickou osobou - podnikatelem, je povinen podat insolvenční návrh bez zbytečného odkladu poté, co se dozvěděl nebo při náležité pečlivosti měl dozvědět o svém úpadku. Tuto povinnost má i tehdy, byl-li pravomocně zastaven výkon rozhodnutí prodejem jeho podniku nebo exekuce podle zvláštního právního předpisu4) proto, že cena majetku náležejícího k podniku nepřevyšuje výši závazků náležejících k podniku; to neplatí, má-li dlužník ještě jiný podnik.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1014"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Pokud osoba, která měla podat návrh (viz přeechozí bod) takto neučinila, odpovídá věřiteli za škodu nebo jinou újmu, kterou způsobí porušením této povinnosti - lze se jí zprostit (pokud by porušení nemělo vliv na konečný rozsah částky vydané k uspokojení věřitelů, nebo že tuto povinnost nesplnila vzhledem ke skutečnostem, které nastaly nezávisle na její vůli a které nemohla odvrátit ani při vynaložení veškerého úsilí, které lze po ní spravedlivě požadovat).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">Soud může návrh posoudit v takzvaném</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">předběžném posouzení insolvenčního návrhu podaného věřitelem</w:t></w:r><w:r><w:t xml:space="preserve">, ten se nemusí okamžitě zveřejňovat v insolvenčním rejstříku (například v případě podezření na šikanózní podání).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Pokud se nejedná o bod výše, uveřejní insolvenční soud oznámení o zahájení insolvenčního řízení v insolvenčním rejstříku.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Navrhovatel (věřitel) musí složit zálohu na náklady insolvenčního řízení.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:bookmarkStart w:id="44" w:name="účinky-spojené-se-zahájením-řízení"/><w:bookmarkEnd w:id="44"/><w:r><w:t xml:space="preserve">Účinky spojené se zahájením řízení</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1015"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">pohledávky nemohou být uplatněny žalobou, lze-li je uplatnit přihláškou,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1015"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">uplatnění ráva ze zajištění podléha podmínkám zákonu,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1015"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">výkon rozhodnutí či exekuci lze nařídit nebo zahájit, nelze ji však provést (mimo případy stanovené zákonem),</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1015"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">nakládání s majetkovou podstatou (podstatné změny) - mimo zákonem stanovené výjimky,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1015"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">soud může rozhodnout o předběžném opatření - například dodatečné omezení dlužníka v dispozici s majetkem náležejícím do majetkové podstaty,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1015"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">a další body.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:bookmarkStart w:id="45" w:name="moratorium"/><w:bookmarkEnd w:id="45"/><w:r><w:t xml:space="preserve">Moratorium</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">Dlužník, který je podnikatelem, může do 7 dnů od podání insolvenčního návrhu, a jde-li o insolvenční návrh věřitele, do 15 dnů od jeho doručení insolvenčním soudem, navrhnout insolvenčnímu soudu vyhlášení moratoria; toto právo nemá právnická osoba v likvidaci.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Moratorium je účinné od okamžiku zveřejnění rozhodnutí o jeho vyhlášení v insolvenčním rejstříku a trvá po dobu uvedenou v návrhu na moratorium, nejdéle však 3 měsíce.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Insolvenční soud může na návrh dlužníka prodloužit moratorium nejdéle o 30 dnů, jestliže dlužník k takovému návrhu připojí ke dni podání návrhu aktualizovaný seznam závazků a písemné prohlášení většiny jeho věřitelů, počítané podle výše jejich pohledávek, že s prodloužením moratoria souhlasí; podpisy věřitelů na tomto prohlášení musí být úředně ověřeny.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Po dobu trvání moratoria nelze vydat rozhodnutí o úpadku.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Není-li dále stanoveno jinak, po dobu trvání moratoria zůstávají zachovány účinky spojené se zahájením insolvenčního řízení.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Účinky moratoria:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1016"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Závazky bezprostředně souvisící se zachováním provozu podniku vzniklé v posledních 30 dnech před vyhlášením moratoria nebo po něm je dlužník po dobu trvání moratoria oprávněn hradit přednostně před dříve splatnými závazky.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1016"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Smlouvy na dodávky energií a surovin, jakož i další smlouvy o dodávkách zboží a služeb, které ke dni vyhlášení moratoria trvaly alespoň po dobu 3 měsíců, nemůže druhý účastník smlouvy po dobu trvání moratoria vypovědět nebo od nich odstoupit pro prodlení dlužníka s placením úhrady za zboží nebo služby, ke kterému došlo před vyhlášením moratoria, nebo pro zhoršení majetkové situace dlužníka, hradí-li dlužník na základě těchto smluv řádně a včas alespoň závazky podle odstavce 1.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1016"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Započtení vzájemných pohledávek dlužníka a věřitele není po dobu trvání moratoria přípustné, ledaže insolvenční soud určí jinak předběžným opatřením. To platí i tehdy, jestliže zákonné podmínky tohoto započtení byly splněny před vyhlášením moratoria.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:bookmarkStart w:id="46" w:name="rozhodnutí-o-úpadku"/><w:bookmarkEnd w:id="46"/><w:r><w:t xml:space="preserve">Rozhodnutí o úpadku</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">Insolvenční soud vydá rozhodnutí o úpadku, je-li osvědčením nebo dokazováním zjištěno, že dlužník je v úpadku nebo že mu úpadek hrozí.</w:t></w:r></w:p>
'@

$wrapper = ('<?xml version="1.0" standalone="yes"?>' + "`n" +
            '<?mso-application progid="Word.Document"?>' + "`n" +
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$endRng = $d.Range($d.Content.End, $d.Content.End)
$endRng.InsertXML($wrapper)

Write-Output "Edit complete. ParaCount=$($d.Paragraphs.Count)"
